# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text flips from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File", "Latest Handback File"
#    and "Latest Handback DateTime" columns populated (+ a hyperlink on the
#    "Latest Target File" cell), and a couple of columns are widened so the
#    longer filenames are readable.

$wb = $excel.ActiveWorkbook

$docUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81ef94ae697e1e8d4be8e68ef2268af8e5437556/e2e/3fec1a29-f8dc-499b-adbd-64b948868ea9.md"
$docName = "3fec1a29-f8dc-499b-adbd-64b948868ea9.md"
$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: widen the per-locale status columns (E, F) and refresh the
# status text that is mirrored there.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Columns.Item(3).ColumnWidth = 29.2

# Latest Target File (I2) -> becomes a hyperlink to the source doc
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $docUrl, "", "", $docName)
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = 15570276

# Latest Handback File (J2) -> the generated zh-cn xliff
$wsZh.Range("J2").Value = "3fec1a29-f8dc-499b-adbd-64b948868ea9.f5863a18e1d03ab84f804234215253f1aa5fb653.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZh.Range("K2").Value = "2016-08-21 17:04:54"

$wsZh.Columns.Item(9).ColumnWidth = 39.2
$wsZh.Columns.Item(10).ColumnWidth = 39.2

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Columns.Item(3).ColumnWidth = 29.2

# Latest Target File (I2) -> becomes a hyperlink to the source doc
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $docUrl, "", "", $docName)
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = 15570276

# Latest Handback File (J2) -> the generated de-de xliff
$wsDe.Range("J2").Value = "3fec1a29-f8dc-499b-adbd-64b948868ea9.f5863a18e1d03ab84f804234215253f1aa5fb653.de-de.xlf"

# Latest Handback DateTime (K2)
$wsDe.Range("K2").Value = "2016-08-21 17:05:02"

$wsDe.Columns.Item(9).ColumnWidth = 39.2
$wsDe.Columns.Item(10).ColumnWidth = 39.2
